$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Fill in the new Status/Description cells for the task log
# (order matters for shared-string table insertion order)
$ws.Range("D6").Value = "DONE"
$ws.Range("E6").Value = "vs code extensions to be installed. "
$ws.Range("E4").Value = "."
$ws.Range("D7").Value = "DONE"
$ws.Range("E7").Value = "installing react js"

# Update selection to match the saved cursor position
$ws.Range("E22").Select()
